$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the "HL7-0190" header label cell (H6) entirely - it is dropped from
# the table and its shared-string entry no longer needed there.
$ws.Range("H6").Clear()

# Add a new footnote cell K6 containing "*ValueはHL7-0190", styled like the
# other plain (non-wrapped) text cells in this area.
$ws.Range("K6").Value = "*ValueはHL7-0190"
$ws.Range("K6").NumberFormat = "@"
$ws.Range("K6").WrapText = $false

# Update the active selection to the new cell.
$ws.Range("K6").Select()
